$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 10770.9
$ws.Range("I9").Value = 15215.286
$ws.Range("K9").Value = 15215.286
$ws.Range("M9").Value = -15046.286

$ws.Range("H17").Value = 401969.66
$ws.Range("J17").Value = 427774.34
$ws.Range("L17").Value = 1283323.02
$ws.Range("N17").Value = -1283659.02

$ws.Range("H33").Value = 263.2353
$ws.Range("I33").Value = 152.16667
$ws.Range("J33").Value = 529.8
$ws.Range("K33").Value = 152.16667
$ws.Range("L33").Value = 529.8
$ws.Range("M33").Value = 76.83332999999999
$ws.Range("N33").Value = -987.8

$ws.Range("H37").Value = 1500
$ws.Range("I37").Value = 1500
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 4500
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -4374
$ws.Range("N37").ClearContents()

$ws.Range("H92").Value = 1001.85
$ws.Range("I92").Value = 952.05554
$ws.Range("J92").Value = 1450
$ws.Range("K92").Value = 952.05554
$ws.Range("L92").Value = 1450
$ws.Range("M92").Value = 295.94446
$ws.Range("N92").Value = -3946

$ws.Range("H106").Value = 3609.125
$ws.Range("I106").Value = 3700.4285
$ws.Range("J106").Value = 2970
$ws.Range("K106").Value = 3700.4285
$ws.Range("L106").Value = 2970
$ws.Range("M106").Value = -3069.4285
$ws.Range("N106").Value = -4232

$ws.Range("H116").Value = 5001
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H129").Value = 9753.299999999999
$ws.Range("I129").Value = 1367
$ws.Range("J129").Value = 43298.5
$ws.Range("K129").Value = 4101
$ws.Range("L129").Value = 129895.5
$ws.Range("M129").Value = 899
$ws.Range("N129").Value = -139895.5

$ws.Range("H132").Value = 2972.739
$ws.Range("J132").Value = 3337.25
$ws.Range("L132").Value = 10011.75
$ws.Range("N132").Value = -15071.75

$ws.Range("H138").Value = 1673.5585
$ws.Range("I138").Value = 1170.0938
$ws.Range("J138").Value = 2031.5778
$ws.Range("K138").Value = 3510.2814
$ws.Range("L138").Value = 6094.7334
$ws.Range("M138").Value = 1629.7186
$ws.Range("N138").Value = -16374.7334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1560.92
$ws.Range("I32").Value = 1593.4681
$ws.Range("K32").Value = 1593.4681
$ws.Range("M32").Value = -1306.4681

$ws.Range("H45").Value = 8108.5713
$ws.Range("I45").Value = 9269.6
$ws.Range("J45").Value = 5206
$ws.Range("K45").Value = 9269.6
$ws.Range("L45").Value = 5206
$ws.Range("M45").Value = -8892.6
$ws.Range("N45").Value = -5960

$ws.Range("H74").Value = 26318022
$ws.Range("I74").Value = 29414128
$ws.Range("K74").Value = 29414128
$ws.Range("M74").Value = -29413254

$ws.Range("H77").Value = 26318022
$ws.Range("I77").Value = 29414128
$ws.Range("K77").Value = 147070640
$ws.Range("M77").Value = -147066272

$ws.Range("H122").Value = 5463.5557
$ws.Range("I122").Value = 3955.3635
$ws.Range("J122").Value = 12099.6
$ws.Range("K122").Value = 11866.0905
$ws.Range("L122").Value = 36298.8
$ws.Range("M122").Value = -9416.0905
$ws.Range("N122").Value = -41198.8

$ws.Range("H132").Value = 3451170.2
$ws.Range("I132").Value = 3451170.2
$ws.Range("K132").Value = 10353510.6
$ws.Range("M132").Value = -10350980.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 779.8182
$ws.Range("I20").Value = 793.8570999999999
$ws.Range("K20").Value = 793.8570999999999
$ws.Range("M20").Value = -546.8570999999999

$ws.Range("H134").Value = 16671508
$ws.Range("I134").Value = 17861998
$ws.Range("J134").Value = 4645.5
$ws.Range("K134").Value = 53585994
$ws.Range("L134").Value = 13936.5
$ws.Range("M134").Value = -53583459
$ws.Range("N134").Value = -19006.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 99999
$ws.Range("I22").Value = 99999
$ws.Range("K22").Value = 99999
$ws.Range("M22").Value = -99649

$ws.Range("H58").Value = 23819716
$ws.Range("I58").Value = 23819716
$ws.Range("K58").Value = 23819716
$ws.Range("M58").Value = -23819513

$ws.Range("H59").Value = 75110.55499999999
$ws.Range("I59").Value = 8666.333000000001
$ws.Range("J59").Value = 108332.664
$ws.Range("K59").Value = 8666.333000000001
$ws.Range("L59").Value = 108332.664
$ws.Range("M59").Value = -7521.333000000001
$ws.Range("N59").Value = -110622.664

$ws.Range("H86").Value = 15035.714
$ws.Range("I86").Value = 10750
$ws.Range("J86").Value = 15750
$ws.Range("K86").Value = 10750
$ws.Range("L86").Value = 15750
$ws.Range("M86").Value = -9627
$ws.Range("N86").Value = -17996

$ws.Range("H89").Value = 15035.714
$ws.Range("I89").Value = 10750
$ws.Range("J89").Value = 15750
$ws.Range("K89").Value = 53750
$ws.Range("L89").Value = 78750
$ws.Range("M89").Value = -48134
$ws.Range("N89").Value = -89982

$ws.Range("H99").Value = 3389.4
$ws.Range("I99").Value = 3600
$ws.Range("K99").Value = 3600
$ws.Range("M99").Value = -2102

$ws.Range("H105").Value = 1901909.8
$ws.Range("J105").Value = 5000
$ws.Range("L105").Value = 5000
$ws.Range("N105").Value = -8494

$ws.Range("H121").Value = 105000
$ws.Range("J121").Value = 105000
$ws.Range("L121").Value = 105000
$ws.Range("N121").Value = -107620

$ws.Range("H126").Value = 3389.4
$ws.Range("I126").Value = 3600
$ws.Range("K126").Value = 10800
$ws.Range("M126").Value = -8330

$ws.Range("H134").Value = 19233306
$ws.Range("I134").Value = 19233306
$ws.Range("K134").Value = 57699918
$ws.Range("M134").Value = -57697383

$ws.Range("H136").Value = 23819716
$ws.Range("I136").Value = 23819716
$ws.Range("K136").Value = 71459148
$ws.Range("M136").Value = -71456598

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H114").Value = 112227.11
$ws.Range("I114").Value = 125755.5
$ws.Range("J114").Value = 4000
$ws.Range("K114").Value = 377266.5
$ws.Range("L114").Value = 12000
$ws.Range("M114").Value = -374012.5
$ws.Range("N114").Value = -18508

$ws.Range("H131").Value = 1616
$ws.Range("I131").Value = 1292
$ws.Range("K131").Value = 3876
$ws.Range("M131").Value = 1164

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7819173
$ws.Range("I132").Value = 10424415
$ws.Range("K132").Value = 31273245
$ws.Range("M132").Value = -31270715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2219.724
$ws.Range("I22").Value = 1897.8572
$ws.Range("K22").Value = 1897.8572
$ws.Range("M22").Value = -1602.8572

$ws.Range("H27").Value = 2219.724
$ws.Range("I27").Value = 1897.8572
$ws.Range("K27").Value = 1897.8572
$ws.Range("M27").Value = -1790.8572

$ws.Range("H93").Value = 1349.75
$ws.Range("I93").Value = 1256.8572
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 1256.8572
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = -8.857199999999921
$ws.Range("N93").Value = -4496

$ws.Range("H136").Value = 1531.0435
$ws.Range("I136").Value = 1405.3889
$ws.Range("J136").Value = 1983.4
$ws.Range("K136").Value = 4216.1667
$ws.Range("L136").Value = 5950.200000000001
$ws.Range("M136").Value = -1666.1667
$ws.Range("N136").Value = -11050.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6666.1113
$ws.Range("I62").Value = 3500
$ws.Range("K62").Value = 3500
$ws.Range("M62").Value = -2876

$ws.Range("H65").Value = 6666.1113
$ws.Range("I65").Value = 3500
$ws.Range("K65").Value = 17500
$ws.Range("M65").Value = -14380

$ws.Range("H113").Value = 1400.091
$ws.Range("I113").Value = 1429.8125
$ws.Range("J113").Value = 1320.8334
$ws.Range("K113").Value = 4289.4375
$ws.Range("L113").Value = 3962.5002
$ws.Range("M113").Value = -2119.4375
$ws.Range("N113").Value = -8302.5002

$ws.Range("H122").Value = 2343.0715
$ws.Range("J122").Value = 2899.5
$ws.Range("L122").Value = 8698.5
$ws.Range("N122").Value = -13598.5

$ws.Range("H126").Value = 3095.6155
$ws.Range("I126").Value = 3232.7222
$ws.Range("J126").Value = 2787.125
$ws.Range("K126").Value = 9698.1666
$ws.Range("L126").Value = 8361.375
$ws.Range("M126").Value = -7228.1666
$ws.Range("N126").Value = -13301.375

$ws.Range("H132").Value = 14713864
$ws.Range("I132").Value = 21742872
$ws.Range("K132").Value = 65228616
$ws.Range("M132").Value = -65226086
